$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-145: only column B (count) values changed
$updatedCounts = @{
    2 = 4923
    3 = 6333
    4 = 7153
    5 = 6136
    6 = 6170
    7 = 5257
    8 = 5822
    9 = 7084
    10 = 6938
    11 = 6035
    12 = 6215
    13 = 5719
    14 = 4692
    15 = 6770
    16 = 5503
    17 = 6759
    18 = 5493
    19 = 5065
    20 = 4797
    21 = 6061
    22 = 4829
    23 = 6076
    24 = 5643
    25 = 4937
    26 = 6403
    27 = 5106
    28 = 6046
    29 = 4817
    30 = 5476
    31 = 5100
    32 = 5511
    33 = 5202
    34 = 5812
    35 = 5057
    36 = 4488
    37 = 4241
    38 = 4591
    39 = 5537
    40 = 4908
    41 = 4781
    42 = 5269
    43 = 4344
    44 = 5203
    45 = 5337
    46 = 4737
    47 = 5479
    48 = 5034
    49 = 5147
    50 = 4983
    51 = 4300
    52 = 4238
    53 = 2648
    54 = 3214
    55 = 5028
    56 = 4606
    57 = 5252
    58 = 4411
    59 = 4319
    60 = 5086
    61 = 4851
    62 = 5366
    63 = 5393
    64 = 6059
    65 = 8486
    66 = 9264
    67 = 10362
    68 = 10264
    69 = 10136
    70 = 12295
    71 = 12764
    72 = 10175
    73 = 10009
    74 = 9911
    75 = 9217
    76 = 8552
    77 = 8559
    78 = 7430
    79 = 7643
    80 = 8285
    81 = 7655
    82 = 10437
    83 = 7503
    84 = 6945
    85 = 6097
    86 = 6597
    87 = 6377
    88 = 6503
    89 = 5444
    90 = 5993
    91 = 8006
    92 = 6750
    93 = 5949
    94 = 6792
    95 = 9145
    96 = 6437
    97 = 7226
    98 = 7065
    99 = 5875
    100 = 5865
    101 = 6659
    102 = 5759
    103 = 5538
    104 = 6293
    105 = 4591
    106 = 4369
    107 = 7036
    108 = 8286
    109 = 7457
    110 = 7029
    111 = 6607
    112 = 7745
    113 = 8526
    114 = 7677
    115 = 7946
    116 = 7126
    117 = 6532
    118 = 6677
    119 = 6885
    120 = 6378
    121 = 6113
    122 = 7287
    123 = 7188
    124 = 6619
    125 = 6623
    126 = 7208
    127 = 5549
    128 = 4887
    129 = 5911
    130 = 4957
    131 = 6486
    132 = 5654
    133 = 5373
    134 = 10244
    135 = 5829
    136 = 5077
    137 = 6497
    138 = 5269
    139 = 4937
    140 = 5053
    141 = 5363
    142 = 5064
    143 = 5921
    144 = 4461
    145 = 4986
}

foreach ($row in $updatedCounts.Keys) {
    $ws.Cells.Item($row, 2).Value = $updatedCounts[$row]
}

# Append new rows 146-180 with ts (A, date-formatted), count (B), time_unit (C = "W")
$newRows = @(
    @(146, 44479, 4559),
    @(147, 44486, 5154),
    @(148, 44493, 5240),
    @(149, 44500, 5334),
    @(150, 44507, 4919),
    @(151, 44514, 4617),
    @(152, 44521, 5053),
    @(153, 44528, 4923),
    @(154, 44535, 5663),
    @(155, 44542, 5331),
    @(156, 44549, 4995),
    @(157, 44556, 4143),
    @(158, 44563, 3656),
    @(159, 44570, 5037),
    @(160, 44577, 4721),
    @(161, 44584, 5410),
    @(162, 44591, 5143),
    @(163, 44598, 5214),
    @(164, 44605, 5662),
    @(165, 44612, 6409),
    @(166, 44619, 5724),
    @(167, 44626, 4871),
    @(168, 44633, 5089),
    @(169, 44640, 5053),
    @(170, 44647, 4693),
    @(171, 44654, 4989),
    @(172, 44661, 5303),
    @(173, 44668, 5177),
    @(174, 44675, 4831),
    @(175, 44682, 5102),
    @(176, 44689, 5080),
    @(177, 44696, 4882),
    @(178, 44703, 5155),
    @(179, 44710, 6182),
    @(180, 44717, 1678),
)

foreach ($item in $newRows) {
    $r = $item[0]
    $dateVal = $item[1]
    $countVal = $item[2]
    $ws.Cells.Item($r, 1).Value = $dateVal
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $countVal
    $ws.Cells.Item($r, 3).Value = "W"
}